# Daily attendance processing - 2025-12-04 03:12:44
# Normalizes the "Recorded By" (column G) cell values on the
# "Session Analysis Results" sheet.
#
# Rule applied to each comma-separated "Recorded By" list:
#   - If the literal token "System" appears anywhere in the list, move it
#     to the end of the list (preserving the relative order of the rest).
#   - Otherwise, if the list has more than one entry, rotate it left by
#     one position (the first recorder moves to the end).
#   - Single-entry lists are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Text

    if ($null -eq $value -or $value -eq "") {
        continue
    }

    $parts = $value -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -le 1) {
        continue
    }

    $systemIndex = -1
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($parts[$i] -ceq "System") {
            $systemIndex = $i
            break
        }
    }

    $newParts = @()
    if ($systemIndex -ge 0) {
        for ($i = 0; $i -lt $parts.Length; $i++) {
            if ($i -ne $systemIndex) {
                $newParts += $parts[$i]
            }
        }
        $newParts += "System"
    } else {
        for ($i = 1; $i -lt $parts.Length; $i++) {
            $newParts += $parts[$i]
        }
        $newParts += $parts[0]
    }

    $newValue = [string]::Join(", ", $newParts)

    if ($newValue -ne $value) {
        $cell.Value = $newValue
    }
}
